# Added Experiment feature + bug fixes
#
# The countries list contained a number of entries that needed to be
# removed (duplicates / no-longer-applicable entries / stray totals row).
# Remove each one by locating its row and deleting the entire row, which
# shifts the remaining countries up and keeps the list contiguous &
# alphabetically ordered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$countriesToRemove = @(
    "Bahamas",
    "Belize",
    "Chad",
    "Djibouti",
    "Dominica",
    "Equatorial Guinea",
    "Grand Total",
    "Guinea",
    "Guinea-Bissau",
    "Liberia",
    "Micronesia, Federated States of...",
    "Monaco",
    "NA",
    "Nauru",
    "Niger",
    "North Korea",
    "Papua New Guinea",
    "Saint Kitts and Nevis",
    "San Marino",
    "Solomon Islands",
    "Swaziland"
)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

foreach ($country in $countriesToRemove) {
    for ($r = $lastRow; $r -ge 1; $r--) {
        $cell = $ws.Cells.Item($r, 1)
        if ($cell.Value2 -eq $country) {
            $ws.Rows.Item($r).Delete()
            break
        }
    }
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
}

$ws.Range("A167:A187").Select()
